{"js": "// Replace the 25 \"two-digit \u00f7 one-digit\" answer strings inside the single\n// table's data rows (rows 0, 4, 8, 12, 16 -- the other rows are blank\n// rows left for students to write their own work), in document order,\n// leaving every other part of the document (fonts, sizes, alignment,\n// the date paragraph, table layout, blank rows) untouched.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// [rowIndex, colIndex, newText] -- 0-based row/col, in document order.\nconst replacements = [\n  [0, 0, \"80\u00f73=26, 2\"],\n  [0, 1, \"10\u00f74=2, 2\"],\n  [0, 2, \"50\u00f78=6, 2\"],\n  [0, 3, \"26\u00f73=8, 2\"],\n  [0, 4, \"89\u00f75=17, 4\"],\n  [4, 0, \"77\u00f72=38, 1\"],\n  [4, 1, \"19\u00f72=9, 1\"],\n  [4, 2, \"93\u00f78=11, 5\"],\n  [4, 3, \"84\u00f76=14, 0\"],\n  [4, 4, \"89\u00f72=44, 1\"],\n  [8, 0, \"96\u00f73=32, 0\"],\n  [8, 1, \"43\u00f73=14, 1\"],\n  [8, 2, \"52\u00f77=7, 3\"],\n  [8, 3, \"53\u00f75=10, 3\"],\n  [8, 4, \"11\u00f73=3, 2\"],\n  [12, 0, \"18\u00f77=2, 4\"],\n  [12, 1, \"56\u00f75=11, 1\"],\n  [12, 2, \"48\u00f76=8, 0\"],\n  [12, 3, \"39\u00f75=7, 4\"],\n  [12, 4, \"87\u00f79=9, 6\"],\n  [16, 0, \"20\u00f77=2, 6\"],\n  [16, 1, \"54\u00f75=10, 4\"],\n  [16, 2, \"77\u00f77=11, 0\"],\n  [16, 3, \"65\u00f76=10, 5\"],\n  [16, 4, \"26\u00f79=2, 8\"],\n];\n\nfor (const [row, col, text] of replacements) {\n  const cell = table.getCell(row, col);\n  const range = cell.body.getRange(\"Whole\");\n  range.insertText(text, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$t.Cell(1, 1).Range.Text = \"80\u00f73=26, 2\"\n$t.Cell(1, 2).Range.Text = \"10\u00f74=2, 2\"\n$t.Cell(1, 3).Range.Text = \"50\u00f78=6, 2\"\n$t.Cell(1, 4).Range.Text = \"26\u00f73=8, 2\"\n$t.Cell(1, 5).Range.Text = \"89\u00f75=17, 4\"\n$t.Cell(5, 1).Range.Text = \"77\u00f72=38, 1\"\n$t.Cell(5, 2).Range.Text = \"19\u00f72=9, 1\"\n$t.Cell(5, 3).Range.Text = \"93\u00f78=11, 5\"\n$t.Cell(5, 4).Range.Text = \"84\u00f76=14, 0\"\n$t.Cell(5, 5).Range.Text = \"89\u00f72=44, 1\"\n$t.Cell(9, 1).Range.Text = \"96\u00f73=32, 0\"\n$t.Cell(9, 2).Range.Text = \"43\u00f73=14, 1\"\n$t.Cell(9, 3).Range.Text = \"52\u00f77=7, 3\"\n$t.Cell(9, 4).Range.Text = \"53\u00f75=10, 3\"\n$t.Cell(9, 5).Range.Text = \"11\u00f73=3, 2\"\n$t.Cell(13, 1).Range.Text = \"18\u00f77=2, 4\"\n$t.Cell(13, 2).Range.Text = \"56\u00f75=11, 1\"\n$t.Cell(13, 3).Range.Text = \"48\u00f76=8, 0\"\n$t.Cell(13, 4).Range.Text = \"39\u00f75=7, 4\"\n$t.Cell(13, 5).Range.Text = \"87\u00f79=9, 6\"\n$t.Cell(17, 1).Range.Text = \"20\u00f77=2, 6\"\n$t.Cell(17, 2).Range.Text = \"54\u00f75=10, 4\"\n$t.Cell(17, 3).Range.Text = \"77\u00f77=11, 0\"\n$t.Cell(17, 4).Range.Text = \"65\u00f76=10, 5\"\n$t.Cell(17, 5).Range.Text = \"26\u00f79=2, 8\"\n"}
